# Apply the edits described by the commit:
# "Added color to indicate the trade to which the data belongs."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet
$ws.Name = "Swap and Legs Schema"

# 2. Update labels to indicate array/list columns
$ws.Range("D12").Value = "leg_type[]"
$ws.Range("E12").Value = "leg_ccy[]"
$ws.Range("D20").Value = "legs[]"

# 3. Color the cells belonging to trade T1 blue, and trade T2 green,
#    across all four diagrams on the sheet.

$blueColor = 15773696    # FF00B0F0 (R=0x00,G=0xB0,B=0xF0) in BGR order for COM
$greenColor = 5287936    # FF00B050 (R=0x00,G=0xB0,B=0x50) in BGR order for COM

$blueRanges = @("B5:C5", "E5:H5", "E6:H6", "B13:E14", "B21:E23", "B31:G31")
$greenRanges = @("B6:C6", "E7:H7", "E8:H8", "B15:E16", "B24:E26", "B32:G32")

foreach ($rng in $blueRanges) {
    $ws.Range($rng).Interior.Color = $blueColor
}

foreach ($rng in $greenRanges) {
    $ws.Range($rng).Interior.Color = $greenColor
}

# 4. Update the sheet view: drop the old scroll/selection position, zoom to 120%
$ws.Activate()
$excel.ActiveWindow.Zoom = 120
